$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5502.609
$ws.Range("I2").Value = 1222.7858
$ws.Range("J2").Value = 12160.111
$ws.Range("K2").Value = 1222.7858
$ws.Range("L2").Value = 12160.111
$ws.Range("M2").Value = -1109.7858
$ws.Range("N2").Value = -12386.111
$ws.Range("H17").Value = 1276.3864
$ws.Range("J17").Value = 1276.3864
$ws.Range("L17").Value = 3829.1592
$ws.Range("N17").Value = -4165.1592
$ws.Range("H21").Value = 5017
$ws.Range("I21").Value = 5017
$ws.Range("K21").Value = 5017
$ws.Range("M21").Value = -4549
$ws.Range("H23").Value = 5017
$ws.Range("I23").Value = 5017
$ws.Range("K23").Value = 5017
$ws.Range("M23").Value = -4783
$ws.Range("H29").Value = 15001.75
$ws.Range("J29").Value = 18002.8
$ws.Range("L29").Value = 54008.39999999999
$ws.Range("N29").Value = -54570.39999999999
$ws.Range("H38").Value = 253.625
$ws.Range("I38").Value = 210.53334
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 631.6000200000001
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -259.6000200000001
$ws.Range("N38").Value = -3444
$ws.Range("H58").Value = 641.8333
$ws.Range("I58").Value = 370.4
$ws.Range("J58").Value = 1999
$ws.Range("K58").Value = 1111.2
$ws.Range("L58").Value = 5997
$ws.Range("M58").Value = -961.1999999999998
$ws.Range("N58").Value = -6297
$ws.Range("H76").Value = 14290665
$ws.Range("I76").Value = 33336534
$ws.Range("J76").Value = 6263
$ws.Range("K76").Value = 33336534
$ws.Range("L76").Value = 6263
$ws.Range("M76").Value = -33336219
$ws.Range("N76").Value = -6893
$ws.Range("H79").Value = 14290665
$ws.Range("I79").Value = 33336534
$ws.Range("J79").Value = 6263
$ws.Range("K79").Value = 33336534
$ws.Range("L79").Value = 6263
$ws.Range("M79").Value = -33335442
$ws.Range("N79").Value = -8447
$ws.Range("H86").Value = 2963.9167
$ws.Range("I86").Value = 1953.6666
$ws.Range("K86").Value = 1953.6666
$ws.Range("M86").Value = -830.6666
$ws.Range("H89").Value = 2963.9167
$ws.Range("I89").Value = 1953.6666
$ws.Range("K89").Value = 9768.333000000001
$ws.Range("M89").Value = -4152.333000000001
$ws.Range("H96").Value = 7053
$ws.Range("I96").Value = 613.75
$ws.Range("J96").Value = 17355.8
$ws.Range("K96").Value = 1841.25
$ws.Range("L96").Value = 52067.39999999999
$ws.Range("M96").Value = -468.25
$ws.Range("N96").Value = -54813.39999999999
$ws.Range("H106").Value = 40485.625
$ws.Range("I106").Value = 45840.715
$ws.Range("K106").Value = 45840.715
$ws.Range("M106").Value = -45209.715
$ws.Range("H138").Value = 2650.756
$ws.Range("I138").Value = 1699.3334
$ws.Range("J138").Value = 3649.75
$ws.Range("K138").Value = 5098.0002
$ws.Range("L138").Value = 10949.25
$ws.Range("M138").Value = 41.9997999999996
$ws.Range("N138").Value = -21229.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2396.25
$ws.Range("J88").Value = 2595
$ws.Range("L88").Value = 2595
$ws.Range("N88").Value = -3407
$ws.Range("H91").Value = 2396.25
$ws.Range("J91").Value = 2595
$ws.Range("L91").Value = 2595
$ws.Range("N91").Value = -5403
$ws.Range("H135").Value = 101738.2
$ws.Range("J135").Value = 101738.2
$ws.Range("L135").Value = 101738.2
$ws.Range("N135").Value = -111878.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 96433.14
$ws.Range("J135").Value = 96433.14
$ws.Range("L135").Value = 96433.14
$ws.Range("N135").Value = -106573.14

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 92541.42999999999
$ws.Range("J125").Value = 92541.42999999999
$ws.Range("L125").Value = 92541.42999999999
$ws.Range("N125").Value = -97461.42999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.666664
$ws.Range("I2").Value = 21.714285
$ws.Range("K2").Value = 21.714285
$ws.Range("M2").Value = 91.285715
$ws.Range("H27").Value = 3306
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 4127.5
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 4127.5
$ws.Range("M27").Value = 146
$ws.Range("N27").Value = -4459.5
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H102").Value = 2372.375
$ws.Range("J102").Value = 1808.6
$ws.Range("L102").Value = 1808.6
$ws.Range("N102").Value = -5052.6
$ws.Range("H122").Value = 30836
$ws.Range("I122").Value = 29669.715
$ws.Range("J122").Value = 39000
$ws.Range("K122").Value = 89009.145
$ws.Range("L122").Value = 117000
$ws.Range("M122").Value = -86559.145
$ws.Range("N122").Value = -121900

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12000
$ws.Range("I20").Value = 12000
$ws.Range("K20").Value = 12000
$ws.Range("M20").Value = -11774
$ws.Range("H22").Value = 3559.2
$ws.Range("I22").Value = 1890
$ws.Range("J22").Value = 5019.75
$ws.Range("K22").Value = 1890
$ws.Range("L22").Value = 5019.75
$ws.Range("M22").Value = -1595
$ws.Range("N22").Value = -5609.75
$ws.Range("H27").Value = 3559.2
$ws.Range("I27").Value = 1890
$ws.Range("J27").Value = 5019.75
$ws.Range("K27").Value = 1890
$ws.Range("L27").Value = 5019.75
$ws.Range("M27").Value = -1783
$ws.Range("N27").Value = -5233.75
$ws.Range("H40").Value = 50009350
$ws.Range("I40").Value = 71438940
$ws.Range("J40").Value = 6966.6665
$ws.Range("K40").Value = 71438940
$ws.Range("L40").Value = 6966.6665
$ws.Range("M40").Value = -71438804
$ws.Range("N40").Value = -7238.6665
$ws.Range("H42").Value = 50000000
$ws.Range("I42").Value = 50000000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 50000000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -49999437
$ws.Range("N42").ClearContents()
$ws.Range("I46").Value = 1800.5
$ws.Range("J46").Value = 7787.2915
$ws.Range("K46").Value = 1800.5
$ws.Range("L46").Value = 7787.2915
$ws.Range("M46").Value = -1612.5
$ws.Range("N46").Value = -8163.2915
$ws.Range("H49").Value = 50000000
$ws.Range("I49").Value = 50000000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 50000000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -49999853
$ws.Range("N49").ClearContents()
$ws.Range("H82").Value = 1524.3636
$ws.Range("I82").Value = 1608.625
$ws.Range("K82").Value = 1608.625
$ws.Range("M82").Value = -1247.625
$ws.Range("H85").Value = 1524.3636
$ws.Range("I85").Value = 1608.625
$ws.Range("K85").Value = 1608.625
$ws.Range("M85").Value = -360.625
$ws.Range("H122").Value = 6164.875
$ws.Range("I122").Value = 6049.2
$ws.Range("K122").Value = 18147.6
$ws.Range("M122").Value = -15697.6
$ws.Range("H125").Value = 210000
$ws.Range("J125").Value = 210000
$ws.Range("L125").Value = 210000
$ws.Range("N125").Value = -219840
$ws.Range("H132").Value = 4280.07
$ws.Range("I132").Value = 3478.76
$ws.Range("K132").Value = 10436.28
$ws.Range("M132").Value = -7906.280000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 774
$ws.Range("I14").Value = 774
$ws.Range("K14").Value = 774
$ws.Range("M14").Value = -606
$ws.Range("H122").Value = 1966.5
$ws.Range("I122").Value = 1759.1666
$ws.Range("J122").Value = 2381.1667
$ws.Range("K122").Value = 5277.4998
$ws.Range("L122").Value = 7143.500100000001
$ws.Range("M122").Value = -2827.4998
$ws.Range("N122").Value = -12043.5001
